$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card6")

# Header M1: "Event " (trailing space) -> "Event" (no trailing space)
$ws.Range("M1").Value = "Event"

# New header N1: "Correction " (keeps a trailing space), styled like the other headers
$ws.Range("N1").Value = "Correction "
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Existing M2:M12 were empty placeholder cells -> now carry "nan" like the rest of the row
$ws.Range("M2:M12").Value = "nan"

# New N2:N12 column cells -> materialize as blank (unstyled) cells
$blankRange = $ws.Range("N2:N12")
$blankRange.Value = ""
$blankRange.Style = "Normal"
